$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
606|44769|Primera|1100|1500|1500|1500|Provincia de Quillota|1500
607|44769|Segunda|1300|1100|1100|1100|Provincia de Quillota|1100
608|44596|Primera|900|1000|1000|1000|Provincia de Quillota|1000
609|44596|Segunda|850|700|700|700|Provincia de Quillota|700
610|44340|Primera|2180|700|800|755|Provincia de Quillota|755
611|44340|Segunda|1300|600|600|600|Provincia de Quillota|600
612|44496|Primera|2480|600|650|624|Provincia de Quillota|624
613|44496|Segunda|1200|500|500|500|Provincia de Quillota|500
614|44399|Primera|1600|700|700|700|Provincia de Quillota|700
615|44399|Segunda|800|550|550|550|Provincia de Quillota|550
616|44425|Primera|1800|700|700|700|Provincia de Quillota|700
617|44425|Segunda|900|600|600|600|Provincia de Quillota|600
618|44377|Primera|1200|600|600|600|Provincia de Quillota|600
619|44377|Segunda|850|700|700|700|Provincia de Quillota|700
620|44512|Primera|2900|600|650|628|Provincia de Quillota|628
621|44512|Segunda|1500|500|500|500|Provincia de Quillota|500
622|44397|Primera|1600|700|700|700|Provincia de Quillota|700
623|44181|Primera|1750|600|650|626|Provincia de Quillota|626
624|44181|Segunda|900|550|550|550|Provincia de Quillota|550
625|44497|Primera|3050|600|650|620|Provincia de Quillota|620
626|44497|Segunda|1300|500|500|500|Provincia de Quillota|500
627|44285|Primera|1600|900|900|900|Provincia de Quillota|900
628|44285|Segunda|1800|700|700|700|Provincia de Quillota|700
629|44362|Primera|1850|700|700|700|Provincia de Quillota|700
630|44557|Primera|2900|700|800|752|Provincia de Quillota|752
631|44557|Segunda|3000|500|600|557|Provincia de Quillota|557
632|44747|Primera|1490|1300|1400|1350|Provincia de Quillota|1350
633|44747|Segunda|850|1000|1000|1000|Provincia de Quillota|1000
634|44357|Primera|1800|800|800|800|Provincia de Santiago|800
635|44357|Segunda|1600|700|700|700|Provincia de Santiago|700
636|44279|Primera|950|1000|1000|1000|Provincia de Quillota|1000
637|44279|Segunda|850|800|800|800|Provincia de Quillota|800
638|44551|Primera|3200|600|700|650|Provincia de Quillota|650
639|44551|Segunda|1500|500|500|500|Provincia de Quillota|500
640|44517|Primera|2600|600|650|625|Provincia de Quillota|625
641|44517|Segunda|1500|500|500|500|Provincia de Quillota|500
642|44757|Primera|2050|1200|1300|1246|Provincia de Quillota|1246
643|44547|Primera|1200|600|600|600|Provincia de Quillota|600
644|44547|Segunda|1600|500|500|500|Provincia de Quillota|500
645|44321|Primera|1300|800|800|800|Provincia de Quillota|800
646|44321|Segunda|1100|700|700|700|Provincia de Quillota|700
647|44438|Primera|3800|600|650|625|Provincia de Quillota|625
648|44438|Segunda|1800|500|500|500|Provincia de Quillota|500
649|44355|Primera|1600|800|800|800|Provincia de Quillota|800
650|44355|Segunda|1200|700|700|700|Provincia de Quillota|700
651|44657|Primera|1600|1100|1200|1147|Provincia de Quillota|1147
652|44657|Segunda|450|900|900|900|Provincia de Quillota|900
653|44391|Primera|1600|600|600|600|Provincia de Quillota|600
654|44391|Segunda|1500|500|500|500|Provincia de Quillota|500
655|44453|Primera|1800|600|650|624|Provincia de Quillota|624
656|44453|Segunda|900|550|550|550|Provincia de Quillota|550
657|44186|Primera|1900|600|650|625|Provincia de Quillota|625
658|44186|Segunda|900|500|500|500|Provincia de Quillota|500
659|44189|Primera|1850|600|700|651|Provincia de Quillota|651
660|44609|Primera|1180|1100|1200|1151|Provincia de Quillota|1151
661|44489|Primera|1750|650|700|674|Provincia de Quillota|674
662|44489|Segunda|900|550|550|550|Provincia de Quillota|550
663|44358|Primera|1800|800|800|800|Provincia de Quillota|800
664|44358|Segunda|450|650|650|650|Provincia de Quillota|650
665|44572|Primera|1750|700|800|749|Provincia de Quillota|749
666|44572|Segunda|900|550|550|550|Provincia de Quillota|550
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split '\|'
    $row = [int]$parts[0]
    $dateVal = [int]$parts[1]
    $calidad = $parts[2]
    $volumen = [double]$parts[3]
    $precioMin = [double]$parts[4]
    $precioMax = [double]$parts[5]
    $precioProm = [double]$parts[6]
    $origen = $parts[7]
    $precioKg = [double]$parts[8]

    $ws.Cells.Item($row, 1).Value = 3                       # A Mercado ID
    $ws.Cells.Item($row, 2).Value = "Femacal de La Calera"  # B Mercado
    $ws.Cells.Item($row, 3).Value = "Coquimbo"              # C Region
    $ws.Cells.Item($row, 4).Value = $dateVal                # D Fecha
    $ws.Cells.Item($row, 5).Value = 5                       # E Codreg
    $ws.Cells.Item($row, 6).Value = 100112006                # F Categoria ID
    $ws.Cells.Item($row, 7).Value = "Repollo"                # G Categoria
    $ws.Cells.Item($row, 8).Value = "Crespo record"          # H Variedad
    $ws.Cells.Item($row, 9).Value = $calidad                 # I Calidad
    $ws.Cells.Item($row, 10).Value = $volumen                # J Volumen
    $ws.Cells.Item($row, 11).Value = $precioMin               # K Precio minimo
    $ws.Cells.Item($row, 12).Value = $precioMax               # L Precio maximo
    $ws.Cells.Item($row, 13).Value = $precioProm              # M Precio promedio ponderado
    $ws.Cells.Item($row, 14).Value = "`$/unidad"              # N Unidad de comercializacion
    $ws.Cells.Item($row, 15).Value = $origen                  # O Origen
    $ws.Cells.Item($row, 16).Value = $precioKg                # P Precio $/Kg
    $ws.Cells.Item($row, 17).Value = 1                        # Q Kg o Unidades
    $ws.Cells.Item($row, 18).Value = "Hortaliza"               # R Clasificacion
}

Write-Host "Updated rows 606-666"
